$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.347.51'
$ws.Range("E2").Value = '  +3.09%  '
$ws.Range("D3").Value = '1.718.04'
$ws.Range("E3").Value = '  +3.13%  '
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '238.67'
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.4721'
$ws.Range("E7").Value = '  -1.67%  '
$ws.Range("D8").Value = '0.2630'
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = '0.06202'
$ws.Range("E9").Value = '  +0.73%  '
$ws.Range("D10").Value = '1.719.07'
$ws.Range("E10").Value = '  +3.21%  '
$ws.Range("D11").Value = '0.07052'
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("D12").Value = '15.30'
$ws.Range("E12").Value = '  +2.84%  '
$ws.Range("D13").Value = '0.5902'
$ws.Range("E13").Value = '  -0.77%  '
$ws.Range("D14").Value = '4.407'
$ws.Range("D15").Value = '76.15'
$ws.Range("E15").Value = '  +2.23%  '
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = '26.350.42'
$ws.Range("E18").Value = '  +3.10%  '
$ws.Range("D19").Value = '0.000006794'
$ws.Range("E19").Value = '  +0.40%  '
$ws.Range("D20").Value = '11.55'
$ws.Range("E20").Value = '  +0.91%  '
$ws.Range("D21").Value = '1.937.09'
$ws.Range("E21").Value = '  +3.01%  '
$ws.Range("D22").Value = '4.544'
$ws.Range("E22").Value = '  +1.83%  '
$ws.Range("D23").Value = '8.748'
$ws.Range("E23").Value = '  +0.65%  '
$ws.Range("D24").Value = '5.320'
$ws.Range("E24").Value = '  -0.22%  '
$ws.Range("D25").Value = '135.36'
$ws.Range("E25").Value = '  +0.33%  '
$ws.Range("D26").Value = '15.24'
$ws.Range("E26").Value = '  +1.08%  '
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").Value = '1.407'
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("B28").Value = 'BitcoinCash'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D28").Value = '108.05'
$ws.Range("E28").Value = '  +2.99%  '
$ws.Range("D29").Value = '1.753'
$ws.Range("E29").Value = '  +3.71%  '
$ws.Range("D30").Value = '3.996'
$ws.Range("E30").Value = '  +0.88%  '
$ws.Range("D31").Value = '3.688'
$ws.Range("E31").Value = '  +0.34%  '
$ws.Range("D32").Value = '0.07736'
$ws.Range("E32").Value = '  +0.85%  '
$ws.Range("D33").Value = '0.04438'
$ws.Range("E33").Value = '  +2.30%  '
$ws.Range("D34").Value = '2.614'
$ws.Range("E34").Value = '  -0.18%  '
$ws.Range("D35").Value = '0.9788'
$ws.Range("E35").Value = '  +2.86%  '
$ws.Range("D36").Value = '0.6185'
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("D37").Value = '0.9294'
$ws.Range("E37").Value = '  +7.91%  '
$ws.Range("D38").Value = '113.68'
$ws.Range("E38").Value = '  +16.29%  '
$ws.Range("D39").Value = '2.409'
$ws.Range("E39").Value = '  -7.40%  '
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("E41").Value = '  +1.53%  '
$ws.Range("D42").Value = '0.01476'
$ws.Range("E42").Value = '  -2.23%  '
$ws.Range("D43").Value = '5.343'
$ws.Range("E43").Value = '  +13.43%  '
$ws.Range("D44").Value = '0.3807'
$ws.Range("E44").Value = '  +0.86%  '
$ws.Range("E45").Value = '  +3.94%  '
$ws.Range("D46").Value = '6.289'
$ws.Range("E46").Value = '  +1.03%  '
$ws.Range("D47").Value = '0.05288'
$ws.Range("E47").Value = '  +0.39%  '
$ws.Range("D48").Value = '30.39'
$ws.Range("E48").Value = '  +3.00%  '
$ws.Range("D49").Value = '7.667'
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '1.217'
$ws.Range("E50").Value = '  +1.50%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '0.3370'
$ws.Range("E51").Value = '  +0.82%  '
